# The "z" column (column E) is removed from Sheet1. Deleting the entire
# column shifts cmvd_Mpc (F->E), RAdeg_gal (G->F) and DEdeg_gal (H->G) one
# position to the left, matching the diff (dimension shrinks to A1:G145).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E:E").Delete()
